$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1)
$ws.Range("G1").Value = "MSE_median"
$ws.Range("H1").Value = "MAE_median"
$ws.Range("I1").Value = "Dir_accuracy"

# Copy the header style from an existing header cell (e.g. F1) to the new headers
$ws.Range("F1").Copy() | Out-Null
$ws.Range("G1:I1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# New data cells (row 2)
$ws.Range("G2").Value = 0.0004614566239965051
$ws.Range("H2").Value = 0.02148153630921317
$ws.Range("I2").Value = 0.5531914893617021

# New data cells (row 3)
$ws.Range("G3").Value = 0.001083194070471167
$ws.Range("H3").Value = 0.03291191380748265
# I3 stays empty (matches F3 which is an empty inlineStr cell) but still needs
# to exist as a blank cell within the sheet's used range, so touch a
# formatting property (without changing the actual formatting) to force
# Excel to materialize the cell record.
$ws.Range("I3").Font.Bold = $false
